$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -0.1510908837489806
$ws.Range("C2").Value = 2.084195974016096
$ws.Range("D2").Value = 20.51875844210171
$ws.Range("E2").Value = 4.529763618788701
$ws.Range("F2").Value = 4.628991593011878
$ws.Range("G2").Value = 23

$ws.Range("B3").Value = -0.2243426211202774
$ws.Range("C3").Value = 2.382314769375817
$ws.Range("D3").Value = 19.09061844583581
$ws.Range("E3").Value = 4.369281227597487
$ws.Range("F3").Value = 4.46620300197608
$ws.Range("G3").Value = 22

$ws.Range("B4").Value = -0.7120868115500602
$ws.Range("C4").Value = 1.933618923903567
$ws.Range("D4").Value = 10.88254014902561
$ws.Range("E4").Value = 3.298869525917267
$ws.Range("F4").Value = 3.300643293046703
$ws.Range("G4").Value = 21

$ws.Range("B5").Value = -0.3051420100610381
$ws.Range("C5").Value = 1.578682252322016
$ws.Range("D5").Value = 10.2251943824045
$ws.Range("E5").Value = 3.197685785439917
$ws.Range("F5").Value = 3.265784782946779
$ws.Range("G5").Value = 20

$ws.Range("B6").Value = -0.2804117665781
$ws.Range("C6").Value = 1.683207935177004
$ws.Range("D6").Value = 11.04924039551323
$ws.Range("E6").Value = 3.324039770446983
$ws.Range("F6").Value = 3.402952828034286
$ws.Range("G6").Value = 19

$ws.Range("B7").Value = -0.3173356620072107
$ws.Range("C7").Value = 1.884213172660337
$ws.Range("D7").Value = 10.94775440033029
$ws.Range("E7").Value = 3.308739095234057
$ws.Range("F7").Value = 3.388969517186058
$ws.Range("G7").Value = 18

$ws.Range("B8").Value = -0.2356007239738259
$ws.Range("C8").Value = 1.940837229566887
$ws.Range("D8").Value = 12.01945225291491
$ws.Range("E8").Value = 3.466908169091721
$ws.Range("F8").Value = 3.565345857874665
$ws.Range("G8").Value = 17

$ws.Range("B9").Value = -0.1538742921933148
$ws.Range("C9").Value = 1.930586848953809
$ws.Range("D9").Value = 10.98523007623265
$ws.Range("E9").Value = 3.314397392623983
$ws.Range("F9").Value = 3.41940389010477
$ws.Range("G9").Value = 16

$ws.Range("B10").Value = -0.1230455590706185
$ws.Range("C10").Value = 1.907399787585887
$ws.Range("D10").Value = 12.46868919605704
$ws.Range("E10").Value = 3.531103113200893
$ws.Range("F10").Value = 3.652819212576486
$ws.Range("G10").Value = 15

$ws.Range("B11").Value = -0.09178185018074302
$ws.Range("C11").Value = 2.13619364861378
$ws.Range("D11").Value = 13.88641155430669
$ws.Range("E11").Value = 3.726447578365579
$ws.Range("F11").Value = 3.865944277603689
$ws.Range("G11").Value = 14
